$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, taken from the target diff.
$updates = @{
    "D2" = "310.17"
    "E2" = "-3.16%"
    "D3" = "49.25"
    "E3" = "1.13%"
    "D4" = "5.170"
    "E4" = "-2.08%"
    "D5" = "0.07777"
    "E5" = "-3.94%"
    "D6" = "4.496"
    "E6" = "-2.27%"
    "D7" = "1.383"
    "E7" = "15.15%"
    "D8" = "1.559"
    "E8" = "-5.84%"
    "D9" = "0.1230"
    "E9" = "-6.66%"
    "D10" = "0.1993"
    "E10" = "2.59%"
    "D11" = "0.09567"
    "E11" = "1.38%"
    "D12" = "0.04686"
    "E12" = "3.98%"
    "E13" = "-0.10%"
    "D14" = "0.001269"
    "E14" = "-4.07%"
    "D15" = "0.04168"
    "E15" = "-2.87%"
    "D16" = "0.005831"
    "E16" = "-1.75%"
    "E17" = "2,018.15%"
    "D18" = "3.334"
    "E18" = "-0.75%"
    "D19" = "2.238"
    "E19" = "-8.12%"
    "D20" = "0.3451"
    "E20" = "1.73%"
    "D21" = "7.923"
    "E21" = "-3.80%"
    "E22" = "-5.12%"
    "D23" = "0.3035"
    "E23" = "-0.80%"
    "D24" = "0.001274"
    "E24" = "-2.76%"
    "D25" = "0.004050"
    "E25" = "-4.62%"
    "D26" = "0.0001351"
    "E26" = "-0.13%"
    "E38" = "-2.54%"
    "D39" = "0.05883"
    "E39" = "5.66%"
    "E40" = "70.03%"
    "D41" = "0.007928"
    "E41" = "3.21%"
    "D42" = "0.1421"
    "E42" = "-1.12%"
    "D43" = "0.008441"
    "E43" = "9.50%"
    "D44" = "0.008333"
    "E44" = "2.86%"
    "D45" = "0.3406"
    "E45" = "6.65%"
    "D46" = "0.00007033"
    "E46" = "0.54%"
    "E47" = "-0.12%"
    "D48" = "0.05190"
    "E48" = "-15.22%"
    "E49" = "-34.58%"
    "E50" = "-0.12%"
    "E51" = "-0.12%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force the cell to stay text-formatted so numeric-looking / percent-looking
    # strings (e.g. "310.17", "-3.16%") are not auto-converted to numbers.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    # Restore the default (unstyled) look so we do not introduce a new style
    # that was not present in the original workbook.
    $range.Style = "Normal"
}
